$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "tv"
$ws.Range("C1").Value = "radio"
$ws.Range("D1").Value = "newspaper"
$ws.Range("E1").Value = "sales"
$ws.Columns("A").Select() | Out-Null
